$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match the repulled data / mean calculation.
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = -3
$ws.Range("F6").Value = 4
$ws.Range("F7").Value = -6
$ws.Range("F8").Value = -1
$ws.Range("F10").Value = 4
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = 3
$ws.Range("F13").Value = 7
$ws.Range("F14").Value = 1
$ws.Range("F15").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = -1
$ws.Range("F18").Value = -3
$ws.Range("F19").Value = -3
$ws.Range("F20").Value = 2
$ws.Range("F21").Value = 5
$ws.Range("F22").Value = -6
$ws.Range("F23").Value = -1
$ws.Range("F24").Value = 1
$ws.Range("F25").Value = 1
$ws.Range("F26").Value = 2
$ws.Range("F27").Value = -6
